# Apply updated crypto price/listing data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'246.97"

# Row 4
$ws.Range("D4").Value = "'5.504"

# Row 5
$ws.Range("D5").Value = "'0.05639"

# Row 6
$ws.Range("D6").Value = "'6.465"

# Row 7
$ws.Range("D7").Value = "'0.8041"

# Row 8
$ws.Range("D8").Value = "'1.053"

# Row 9
$ws.Range("D9").Value = "'0.1453"

# Row 10
$ws.Range("D10").Value = "'0.07327"

# Row 11
$ws.Range("D11").Value = "'0.03175"

# Row 12
$ws.Range("B12").Value = "ProBitToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D12").Value = "'0.1313"
$ws.Range("E12").Value = "11ProBitTokenPROBBestin24h"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02926"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09253"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001664"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Row 16
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.205"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04724"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005840"
$ws.Range("E18").Value = "17OneONE"

# Row 19
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006280"
$ws.Range("E19").Value = "18TigerCashTCH"

# Row 20
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001057"
$ws.Range("E20").Value = "19BitKanKAN"

# Row 21
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.004114"
$ws.Range("E21").Value = "20HotbitTokenHTB"

# Row 22
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001503"
$ws.Range("E22").Value = "21NitroExNTX"

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.976"
$ws.Range("E23").Value = "22LEOLEO"

# Row 24
$ws.Range("B24").Value = "GateToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D24").Value = "'3.380"
$ws.Range("E24").Value = "23GateTokenGT"

# Row 25
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.138"
$ws.Range("E25").Value = "24BTSETokenBTSE"

# Row 26
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3274"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"

# Row 27
$ws.Range("D27").Value = "'0.0003005"

# Row 40
$ws.Range("D40").Value = "'0.04156"

# Row 41
$ws.Range("D41").Value = "'0.006923"

# Row 42
$ws.Range("D42").Value = "'0.003506"

# Row 43
$ws.Range("D43").Value = "'0.1038"

# Row 44
$ws.Range("D44").Value = "'0.009838"

# Row 45
$ws.Range("D45").Value = "'0.00005640"

# Row 46
$ws.Range("D46").Value = "'0.00000000751"

# Row 47
$ws.Range("D47").Value = "'0.6812"

# Row 48
$ws.Range("D48").Value = "'0.02133"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

# Row 49
$ws.Range("D49").Value = "'0.00002104"
